$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.060212565978397
$ws.Range("D2").Value = 1.059029976307232
$ws.Range("E2").Value = 1.065305517444204
$ws.Range("F2").Value = 1.074369071146583
$ws.Range("I2").Value = 1.048266099844684
$ws.Range("J2").Value = 1.065194583162611
$ws.Range("K2").Value = 1.061760593960951
$ws.Range("L2").Value = 1.068019114961758
$ws.Range("M2").Value = 1.077058466466709
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.061496188969643
$ws.Range("D3").Value = 1.060029360536541
$ws.Range("E3").Value = 1.066488710491749
$ws.Range("F3").Value = 1.075724833772661
$ws.Range("I3").Value = 1.048655696046671
$ws.Range("J3").Value = 1.06613017084674
$ws.Range("K3").Value = 1.062573658481804
$ws.Range("L3").Value = 1.06901676174861
$ws.Range("M3").Value = 1.078230018876174
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.062326228587303
$ws.Range("D4").Value = 1.060675445267673
$ws.Range("E4").Value = 1.067254098604846
$ws.Range("F4").Value = 1.076602188036039
$ws.Range("I4").Value = 1.04890627515856
$ws.Range("J4").Value = 1.066734492371149
$ws.Range("K4").Value = 1.063098569188279
$ws.Range("L4").Value = 1.069661512171499
$ws.Range("M4").Value = 1.078987631889343
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.062675048634612
$ws.Range("D5").Value = 1.060946921260194
$ws.Range("E5").Value = 1.067575818131303
$ws.Range("F5").Value = 1.07697105099508
$ws.Range("I5").Value = 1.049011256702166
$ws.Range("J5").Value = 1.066988295935871
$ws.Range("K5").Value = 1.063318956924476
$ws.Range("L5").Value = 1.069932377336073
$ws.Range("M5").Value = 1.079306024841854
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.062733609604773
$ws.Range("D6").Value = 1.06099249522939
$ws.Range("E6").Value = 1.067629833411248
$ws.Range("F6").Value = 1.077032986153258
$ws.Range("I6").Value = 1.049028862369758
$ws.Range("J6").Value = 1.067030895883025
$ws.Range("K6").Value = 1.063355944347982
$ws.Range("L6").Value = 1.069977845831114
$ws.Range("M6").Value = 1.079359478184986
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.062330890039277
$ws.Range("D7").Value = 1.060679073283402
$ws.Range("E7").Value = 1.067258397629082
$ws.Range("F7").Value = 1.076607116706147
$ws.Range("I7").Value = 1.04890767934799
$ws.Range("J7").Value = 1.066737884699153
$ws.Range("K7").Value = 1.063101515136658
$ws.Range("L7").Value = 1.069665132220821
$ws.Range("M7").Value = 1.078991886691065
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.060646488206924
$ws.Range("D8").Value = 1.059367844774464
$ws.Range("E8").Value = 1.06570542936734
$ws.Range("F8").Value = 1.074827241541854
$ws.Range("I8").Value = 1.048398080105304
$ws.Range("J8").Value = 1.065510991421919
$ws.Range("K8").Value = 1.062035621669734
$ws.Range("L8").Value = 1.068356439991605
$ws.Range("M8").Value = 1.077454494788274
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.05767400184968
$ws.Range("D9").Value = 1.057052745184799
$ws.Range("E9").Value = 1.06296713660516
$ws.Range("F9").Value = 1.071691399359271
$ws.Range("I9").Value = 1.047488450379145
$ws.Range("J9").Value = 1.063340789633618
$ws.Range("K9").Value = 1.060148141737401
$ws.Range("L9").Value = 1.06604418278812
$ws.Range("M9").Value = 1.074741780369689
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.055689217650336
$ws.Range("D10").Value = 1.055506177875344
$ws.Range("E10").Value = 1.061140275568745
$ws.Range("F10").Value = 1.069601011056389
$ws.Range("I10").Value = 1.046874136700831
$ws.Range("J10").Value = 1.06188830382356
$ws.Range("K10").Value = 1.058883500857672
$ws.Range("L10").Value = 1.064498391807705
$ws.Range("M10").Value = 1.072930701826595
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.054828997794873
$ws.Range("D11").Value = 1.054835720577301
$ws.Range("E11").Value = 1.060348875237239
$ws.Range("F11").Value = 1.068695850170939
$ws.Range("I11").Value = 1.046606245906017
$ws.Range("J11").Value = 1.061257983178738
$ws.Range("K11").Value = 1.058334373306913
$ws.Range("L11").Value = 1.06382799995354
$ws.Range("M11").Value = 1.072145831480218
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054509350421027
$ws.Range("D12").Value = 1.054586562907234
$ws.Range("E12").Value = 1.060054857256318
$ws.Range("F12").Value = 1.068359628171486
$ws.Range("I12").Value = 1.046506454288682
$ws.Range("J12").Value = 1.061023643320154
$ws.Range("K12").Value = 1.05813017083365
$ws.Range("L12").Value = 1.063578825462843
$ws.Range("M12").Value = 1.071854193119243
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054577921516277
$ws.Range("D13").Value = 1.054640013539115
$ws.Range("E13").Value = 1.060117927722948
$ws.Range("F13").Value = 1.068431749213378
$ws.Range("I13").Value = 1.046527872848666
$ws.Range("J13").Value = 1.06107391959319
$ws.Range("K13").Value = 1.058173983494062
$ws.Range("L13").Value = 1.063632281571344
$ws.Range("M13").Value = 1.071916755230633
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054802578183818
$ws.Range("D14").Value = 1.054815127587263
$ws.Range("E14").Value = 1.060324572792232
$ws.Range("F14").Value = 1.068668058077305
$ws.Range("I14").Value = 1.046598002921999
$ws.Range("J14").Value = 1.061238616886746
$ws.Range("K14").Value = 1.058317498613476
$ws.Range("L14").Value = 1.063807406414777
$ws.Range("M14").Value = 1.072121726673243
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.054940980098828
$ws.Range("D15").Value = 1.054923005165436
$ws.Range("E15").Value = 1.060451886042767
$ws.Range("F15").Value = 1.068813655004493
$ws.Range("I15").Value = 1.046641174565337
$ws.Range("J15").Value = 1.06134006432435
$ws.Range("K15").Value = 1.058405892199252
$ws.Range("L15").Value = 1.063915285180548
$ws.Range("M15").Value = 1.072248002656694
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.055746290518435
$ws.Range("D16").Value = 1.055550657218012
$ws.Range("E16").Value = 1.061192790419112
$ws.Range("F16").Value = 1.069661083107728
$ws.Range("I16").Value = 1.046891875825856
$ws.Range("J16").Value = 1.06193010672708
$ws.Range("K16").Value = 1.05891991222546
$ws.Range("L16").Value = 1.064542861064381
$ws.Range("M16").Value = 1.072982776878892
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.056251224787447
$ws.Range("D17").Value = 1.05594415541545
$ws.Range("E17").Value = 1.061657442376085
$ws.Range("F17").Value = 1.070192647433422
$ws.Range("I17").Value = 1.047048627498005
$ws.Range("J17").Value = 1.06229985213106
$ws.Range("K17").Value = 1.059241932290571
$ws.Range("L17").Value = 1.064936238682391
$ws.Range("M17").Value = 1.073443502006855
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.056545667946984
$ws.Range("D18").Value = 1.056173600820468
$ws.Range("E18").Value = 1.061928431714366
$ws.Range("F18").Value = 1.070502699292043
$ws.Range("I18").Value = 1.047139875928965
$ws.Range("J18").Value = 1.062515384947563
$ws.Range("K18").Value = 1.059429613704123
$ws.Range("L18").Value = 1.065165587666517
$ws.Range("M18").Value = 1.073712171582963
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056646052605709
$ws.Range("D19").Value = 1.056251823047091
$ws.Range("E19").Value = 1.062020826348939
$ws.Range("F19").Value = 1.070608418966554
$ws.Range("I19").Value = 1.047170958424934
$ws.Range("J19").Value = 1.062588853517665
$ws.Range("K19").Value = 1.059493583274224
$ws.Range("L19").Value = 1.065243772604865
$ws.Range("M19").Value = 1.073803770274621
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.056197058057251
$ws.Range("D20").Value = 1.055901944603071
$ws.Range("E20").Value = 1.061607593187574
$ws.Range("F20").Value = 1.070135615703382
$ws.Range("I20").Value = 1.047031828379147
$ws.Range("J20").Value = 1.062260195769278
$ws.Range("K20").Value = 1.059207397852859
$ws.Range("L20").Value = 1.064894043508917
$ws.Range("M20").Value = 1.073394077163324
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054736425833378
$ws.Range("D21").Value = 1.054763564194969
$ws.Range("E21").Value = 1.060263722557756
$ws.Range("F21").Value = 1.068598471172181
$ws.Range("I21").Value = 1.046577359242218
$ws.Range("J21").Value = 1.061190123476784
$ws.Range("K21").Value = 1.058275243411534
$ws.Range("L21").Value = 1.063755840982589
$ws.Range("M21").Value = 1.072061370549332
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.053817352356639
$ws.Range("D22").Value = 1.054047123963296
$ws.Range("E22").Value = 1.059418448459404
$ws.Range("F22").Value = 1.067631976998143
$ws.Range("I22").Value = 1.046289966629255
$ws.Range("J22").Value = 1.060516106145079
$ws.Range("K22").Value = 1.057687817230076
$ws.Range("L22").Value = 1.063039274403961
$ws.Range("M22").Value = 1.071222849984825
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.054304639788188
$ws.Range("D23").Value = 1.054426989152226
$ws.Range("E23").Value = 1.059866576455679
$ws.Range("F23").Value = 1.068144337932157
$ws.Range("I23").Value = 1.046442475697203
$ws.Range("J23").Value = 1.060873532059885
$ws.Range("K23").Value = 1.057999351075967
$ws.Range("L23").Value = 1.063419229302021
$ws.Range("M23").Value = 1.071667423074169
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.056221533907146
$ws.Range("D24").Value = 1.055921018085319
$ws.Range("E24").Value = 1.061630117996018
$ws.Range("F24").Value = 1.070161385889107
$ws.Range("I24").Value = 1.047039419740694
$ws.Range("J24").Value = 1.062278115184366
$ws.Range("K24").Value = 1.059223002934276
$ws.Range("L24").Value = 1.064913110005173
$ws.Range("M24").Value = 1.073416410319069
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.058442997324307
$ws.Range("D25").Value = 1.057651804118409
$ws.Range("E25").Value = 1.063675275024913
$ws.Range("F25").Value = 1.072502045977658
$ws.Range("I25").Value = 1.047724998360827
$ws.Range("J25").Value = 1.063902831083202
$ws.Range("K25").Value = 1.060637206640578
$ws.Range("L25").Value = 1.066642702244386
$ws.Range("M25").Value = 1.075443529022866
